# Apply the "allow-gold as 10000" edit:
#  - set D2:D11 on the "jobs" sheet to 10000
#  - select D2:D11 (active cell D2) on "jobs" and make it the active/selected sheet
#  - un-select/tab the "notices" sheet (it was previously the active tab)
#  - adjust the workbook window size / drop the stored activeTab

$wb = $excel.ActiveWorkbook

$jobs    = $wb.Worksheets.Item("jobs")
$skills  = $wb.Worksheets.Item("skills")
$notices = $wb.Worksheets.Item("notices")

# Update the "Allowed Gold" column values to 10000 for all data rows.
$jobs.Range("D2:D11").Value = 10000

# Activate the jobs sheet and select D2:D11 so it becomes the active/selected range & tab.
$jobs.Activate()
$jobs.Range("D2:D11").Select()

# Resize the workbook window to match the new view geometry.
$win = $excel.ActiveWindow
$win.Width = 20700
$win.Height = 9870

$wb.Save()
